{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows;\nrows.items.forEach((r) => r.cells.load(\"items\"));\nawait context.sync();\n\n// --- Bottom of the table first (rows 24, 25, 26 / index 23, 24, 25) so these\n//     edits don't shift the indices used for the row1/row3 edits below. ---\n\n// Row 26 (index 25): ten tab-separated \"0\" values collapse into a single \"63\"\nrows.items[25].cells.items[0].value = \"63\";\n\n// Rows 24+25 merge: row24's tab-separated values are replaced by a single\n// \"100\", and the (empty) row25 is removed entirely.\nrows.items[23].cells.items[0].value = \"100\";\nrows.items[24].delete();\nawait context.sync();\n\n// --- Row 3 (index 2): \"63\" -> \"0.00000\", with 9 new rows inserted right after it ---\nrows.items[2].cells.items[0].value = \"0.00000\";\nrows.items[2].insertRows(\"After\", 9, [\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.0\"],\n]);\n\n// --- Row 1 (index 0): \"100\" -> \"0M\", with 12 new rows inserted right after it ---\nrows.items[0].cells.items[0].value = \"0M\";\nrows.items[0].insertRows(\"After\", 12, [\n  [\"0M\"],\n  [\"0M\"],\n  [\"10\"],\n  [\"0.00003\"],\n  [\"0.00005\"],\n  [\"0.00003\"],\n  [\"0.00001\"],\n  [\"0.00003\"],\n  [\"0.00003\"],\n  [\"0.00004\"],\n  [\"0.00035\"],\n  [\"100.0\"],\n]);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- Bottom of the table first (indices 24, 25, 26) so earlier edits don't shift these rows ---\n\n# Row 26: ten tab-separated \"0\" values collapse into a single \"63\"\n$t.Rows.Item(26).Cells.Item(1).Range.Text = \"63\"\n\n# Rows 24+25 merge: row24's tab-separated values are replaced by a single \"100\",\n# and the (empty) row25 is removed entirely.\n$t.Rows.Item(24).Cells.Item(1).Range.Text = \"100\"\n$t.Rows.Item(25).Delete()\n\n# --- Row 3: \"63\" -> \"0.00000\", with 9 new rows inserted right after it ---\n$t.Rows.Item(3).Cells.Item(1).Range.Text = \"0.00000\"\n\n$row3Values = @(\"0.00000\", \"0.00000\", \"0.00000\", \"0.00000\", \"0.00000\", \"0.00000\", \"0.00000\", \"0.00000\", \"0.0\")\n$afterRow = $t.Rows.Item(4)\nfor ($i = $row3Values.Count - 1; $i -ge 0; $i--) {\n    $newRow = $t.Rows.Add($afterRow)\n    $newRow.Cells.Item(1).Range.Text = $row3Values[$i]\n}\n\n# --- Row 1: \"100\" -> \"0M\", with 12 new rows inserted right after it ---\n$t.Rows.Item(1).Cells.Item(1).Range.Text = \"0M\"\n\n$row1Values = @(\"0M\", \"0M\", \"10\", \"0.00003\", \"0.00005\", \"0.00003\", \"0.00001\", \"0.00003\", \"0.00003\", \"0.00004\", \"0.00035\", \"100.0\")\n$afterRow = $t.Rows.Item(2)\nfor ($i = $row1Values.Count - 1; $i -ge 0; $i--) {\n    $newRow = $t.Rows.Add($afterRow)\n    $newRow.Cells.Item(1).Range.Text = $row1Values[$i]\n}\n"}
